$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldPrefix = "/localhost:8000/media/"
$newPrefix = "https://raw.githubusercontent.com/xvang3/HLLA/sqlite-testing/hmoob_lus/media/"

# Update all Male/Female Audio URL cells in columns D and E (rows 2-35)
foreach ($col in @("D", "E")) {
    for ($row = 2; $row -le 35; $row++) {
        $cell = $ws.Range("$col$row")
        $val = $cell.Value2
        if ($val -and $val.ToString().StartsWith($oldPrefix)) {
            $cell.Value2 = $val.ToString().Replace($oldPrefix, $newPrefix)
        }
    }
}

# Widen columns D and E so the longer URLs fit (matches the post-edit bestFit
# widths of 125.1796875 and 127 character-units respectively)
$ws.Columns.Item(4).ColumnWidth = 124.25
$ws.Columns.Item(5).ColumnWidth = 126.1

# Update the sheet view selection to a full-column selection (A1:XFD1048576)
$ws.Cells.EntireColumn.Select() | Out-Null
